$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "UVic Support Team" phishing example -> append "Slovakia" as a new
# line to the raw text, and flag the RESULT column as "Y" (True).
$a4 = $ws.Cells.Item(4, 1)
$a4.Value2 = $a4.Value2 + "`nSlovakia"
$ws.Cells.Item(4, 2).Value2 = "Y"
# Re-fit the row so the extra line doesn't leave a stray explicit row height.
$ws.Rows.Item(4).AutoFit()

# Remove the "Internal Revenue Service (IRS)" phishing example entirely.
$ws.Rows.Item(5).Delete()

# Add two more countries with "Y" (True) results.
$ws.Cells.Item(11, 1).Value2 = "United Kingdom"
$ws.Cells.Item(11, 2).Value2 = "Y"
$ws.Cells.Item(12, 1).Value2 = "Nigeria"
$ws.Cells.Item(12, 2).Value2 = "Y"
